$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.396.99'
$ws.Range('E2').Value = '  +5.14%  '
$ws.Range('D3').Value = '2.467.08'
$ws.Range('E3').Value = '  +6.50%  '
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '567.85'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +4.14%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '143.66'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +10.57%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('E8').Value = '  +2.34%  '
$ws.Range('D9').Value = '2.465.56'
$ws.Range('E9').Value = '  +6.45%  '
$ws.Range('E10').Value = '  +5.02%  '
$ws.Range('E11').Value = '  +3.54%  '
$ws.Range('E12').Value = '  +1.29%  '
$ws.Range('E13').Value = '  +5.43%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '26.48'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +12.98%  '
$ws.Range('D15').Value = '2.906.86'
$ws.Range('E15').Value = '  +6.21%  '
$ws.Range('D16').Value = '63.240.39'
$ws.Range('E16').Value = '  +4.82%  '
$ws.Range('E17').Value = '  +7.00%  '
$ws.Range('D18').Value = '2.467.15'
$ws.Range('E18').Value = '  +6.21%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.27'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +6.56%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '341.80'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +9.15%  '
$ws.Range('E21').Value = '  +5.70%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.81'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +3.60%  '
$ws.Range('E23').Value = '  -0.08%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '65.67'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.57%  '
$ws.Range('E25').Value = '  +2.53%  '
$ws.Range('E26').Value = '  +0.03%  '
$ws.Range('E27').Value = '  +9.24%  '
$ws.Range('E28').Value = '  +3.28%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.35'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +9.58%  '
$ws.Range('B30').Value = 'Aptos'
$ws.Range('C30').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.88'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +15.19%  '
$ws.Range('B31').Value = 'PEPE'
$ws.Range('C31').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D31').Value = '0.0₃0816'
$ws.Range('E31').Value = '  +12.53%  '
$ws.Range('E32').Value = '  +7.37%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '176.01'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.67%  '
$ws.Range('E34').Value = '  +11.11%  '
$ws.Range('E35').Value = '  +4.76%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '18.93'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +5.35%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '372.22'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +16.87%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.46'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +8.62%  '
$ws.Range('E39').Value = '  +0.01%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.999'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.05%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.72'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +12.76%  '
$ws.Range('E42').Value = '  +6.36%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '151.28'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +10.38%  '
$ws.Range('E44').Value = '  +6.14%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '20.62'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +8.07%  '
$ws.Range('E46').Value = '  +6.28%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0963'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.45%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0521'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +5.09%  '
$ws.Range('E49').Value = '  +9.44%  '
$ws.Range('E50').Value = '  +4.87%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '18.04'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +6.93%  '
